# "Saving G_12 to G_96 abstraction results."
#
# The author pasted two more blocks of terminal output (tmux capture-pane
# text, two columns separated by a box-drawing "│") from the solver/results
# session into the "AISG Abstract" sheet, continuing the existing log below
# row 125: one block for the G_12->G_64 abstraction run, one for the
# G_12->G_96 abstraction run. Each block leaves a 3-row gap before it
# (rows 126-128, 147-149) and starts with a long dashed separator line that
# Excel's paste parser mis-detects as a leading-minus array formula (it
# evaluates to an error because of the non-formula "│" token).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AISG Abstract")
$ws.Activate()

# --- G_12 -> G_64 block (rows 129-146) ---------------------------------
$ws.Range("A129").FormulaArray = "=" + '---------------------------------------------------------------------------------------------                          │'
$ws.Range("A130").Value = 'Hi: 17227064385339908, Lo: 1, Resolution: 1, Max Time: 5:00:00                                                         │[jtsoundy@hopper:~/Projects/h-generator/solver/results]$ cp abstract/G_12/5H/config_G_48.5H.abstract.from_G12 ~/ProjeTimeout: 5:00:00, Update Time: 0:00:00.010000, SAT Update Time: 0:00:00.010000                                         │cts/temp/abstract/G_12/5H/'
$ws.Range("A131").Value = 'Stride discount: 0.5, Stutter discount: 0.75                                                                           │'
$ws.Range("A132").Value = 'Reserved finger: pinky                                                                                                 │[jtsoundy@hopper:~/Projects/h-generator/solver/results]$ git fetch'
$ws.Range("A133").Value = '---------------------------------------------------------------------------------------------                          │remote: Enumerating objects: 16, done.'
$ws.Range("A134").Value = 'N-Grams: 64, Setup Time: 0:00:06.712160, Current Time: 2022-04-29 11:17:05.072082                                      │remote: Counting objects: 100% (16/16), done.'
$ws.Range("A135").Value = '---------------------------------------------------------------------------------------------                          │remote: Compressing objects: 100% (3/3), done.'
$ws.Range("A136").Value = 'Cost Constraint         - Actual Cost             - Result  - Time:This Run  - Time:All Runs                           │remote: Total 10 (delta 8), reused 9 (delta 7), pack-reused 0'
$ws.Range("A137").Value = '17,227,064,040,798,620  - 13,939,817,899,413,599  - sat     - 0:00:02.330459 - 0:00:02.330491                          │Unpacking objects: 100% (10/10), 22.04 KiB | 663.00 KiB/s, done.'
$ws.Range("A138").Value = '13,939,817,727,142,956  - 13,864,577,058,689,059  - sat     - 0:00:02.968555 - 0:00:05.374564                          │From https://github.com/HiDefender/h-generator'
$ws.Range("A139").Value = '13,864,576,886,418,416  - 13,848,833,915,267,158  - sat     - 0:00:03.821239 - 0:00:09.271975                          │   d6be64a..f75071f  z3-twiddler-model -> origin/z3-twiddler-model'
$ws.Range("A140").Value = '13,848,833,742,996,514  - 13,770,128,533,343,768  - sat     - 0:00:00.815063 - 0:00:10.163481                          │'
$ws.Range("A141").Value = '13,770,128,361,073,124  - 13,682,419,186,271,063  - sat     - 0:00:01.135667 - 0:00:11.375883                          │[jtsoundy@hopper:~/Projects/h-generator/solver/results]$ git reset --hard origin/z3-twiddler-model'
$ws.Range("A142").Value = '13,682,419,014,000,420  - 13,682,419,186,271,063  - unknown - 4:59:48.143567 - 4:59:59.596067                          │HEAD is now at f75071f Setup for G_12 to G_64'
$ws.Range("A143").FormulaArray = "=" + '---------------------------------------------------------------------------------------------                          │'
$ws.Range("A144").Value = 'Sat: 13682419186271063, Unknown: 1.368241901400042e+16, Unsat: 0                                                       │[jtsoundy@hopper:~/Projects/h-generator/solver/results]$ nano ../lib/parameters.py'
$ws.Range("A145").Value = 'Total Time: 5:00:07.345937                                                                                             │'
$ws.Range("A146").Value = '---------------------------------------------------------------------------------------------'

# --- G_12 -> G_96 block (rows 150-169) ---------------------------------
$ws.Range("A150").FormulaArray = "=" + '---------------------------------------------------------------------------------------------                          │config_G_24.5H.abstract.from_G12  config_G_48.5H.abstract.from_G12'
$ws.Range("A151").Value = 'Hi: 17227064385339908, Lo: 1, Resolution: 1, Max Time: 5:00:00                                                         │'
$ws.Range("A152").Value = 'Timeout: 5:00:00, Update Time: 0:00:00.010000, SAT Update Time: 0:00:00.010000                                         │[jtsoundy@hopper:~/Projects/h-generator/solver/results]$ cp abstract/G_12/5H/config_G_64.5H.abstract.from_G12 ~/Proje'
$ws.Range("A153").Value = 'Stride discount: 0.5, Stutter discount: 0.75                                                                           │cts/temp/abstract/G_12/5H/'
$ws.Range("A154").Value = 'Reserved finger: pinky                                                                                                 │'
$ws.Range("A155").Value = '---------------------------------------------------------------------------------------------                          │[jtsoundy@hopper:~/Projects/h-generator/solver/results]$ git fetch'
$ws.Range("A156").Value = 'N-Grams: 96, Setup Time: 0:00:11.346711, Current Time: 2022-04-29 17:27:31.955427                                      │remote: Enumerating objects: 11, done.'
$ws.Range("A157").Value = '---------------------------------------------------------------------------------------------                          │remote: Counting objects: 100% (11/11), done.'
$ws.Range("A158").Value = 'Cost Constraint         - Actual Cost             - Result  - Time:This Run  - Time:All Runs                           │remote: Compressing objects: 100% (1/1), done.'
$ws.Range("A159").Value = '17,227,064,040,798,620  - 13,848,697,624,521,471  - sat     - 0:00:08.733732 - 0:00:08.733756                          │remote: Total 6 (delta 5), reused 6 (delta 5), pack-reused 0'
$ws.Range("A160").Value = '13,848,697,452,250,828  - 13,845,839,182,776,341  - sat     - 0:00:18.630233 - 0:00:27.491627                          │Unpacking objects: 100% (6/6), 1.06 KiB | 361.00 KiB/s, done.'
$ws.Range("A161").Value = '13,845,839,010,505,696  - 13,766,611,927,037,079  - sat     - 0:00:19.331587 - 0:00:46.953769                          │From https://github.com/HiDefender/h-generator'
$ws.Range("A162").Value = '13,766,611,754,766,436  - 13,726,944,792,577,171  - sat     - 0:00:10.495408 - 0:00:57.581177                          │   f75071f..23c0d1d  z3-twiddler-model -> origin/z3-twiddler-model'
$ws.Range("A163").Value = '13,726,944,620,306,528  - 13,716,371,491,029,603  - sat     - 0:07:45.784454 - 0:08:43.496720                          │'
$ws.Range("A164").Value = '13,716,371,318,758,960  - 13,706,703,645,130,208  - sat     - 2:13:34.357549 - 2:22:17.985970                          │[jtsoundy@hopper:~/Projects/h-generator/solver/results]$ git reset --hard origin/z3-twiddler-model'
$ws.Range("A165").Value = '13,706,703,472,859,564  - 13,706,703,645,130,208  - unknown - 2:37:41.062116 - 4:59:59.179632                          │HEAD is now at 23c0d1d Setup for G_12 to G_96'
$ws.Range("A166").FormulaArray = "=" + '---------------------------------------------------------------------------------------------                          │'
$ws.Range("A167").Value = 'Sat: 13706703645130208, Unknown: 1.3706703472859564e+16, Unsat: 0                                                      │[jtsoundy@hopper:~/Projects/h-generator/solver/results]$ nano ../lib/buttons.py'
$ws.Range("A168").Value = 'Total Time: 5:00:13.420140                                                                                             │'
$ws.Range("A169").Value = '---------------------------------------------------------------------------------------------'

# Scroll/select to match where the author left the view on this sheet.
$excel.ActiveWindow.ScrollRow = 127
$ws.Range("I150").Select() | Out-Null

# The "AISG Time Results" tab's view also scrolled down a few rows while
# the workbook was open (selection itself is unchanged, still M45).
$ws3 = $wb.Worksheets.Item("AISG Time Results")
$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 28
$ws3.Range("M45").Select() | Out-Null

# Re-activate "AISG Abstract", which is the tab that was active/saved.
$ws.Activate()

# Recalculate so the volatile RAND()-based cost-jitter cells on
# "Generate Cost Function" (V32/V33) pick up fresh values, as happens
# whenever the workbook is recalculated and saved.
$excel.CalculateFull()
